# Sprint 39 - Day 9 Test Case Summary numbers were filled in after the
# "created test cases for signup (positive and negative) and step1,2 for
# fibmashvpn" work was completed: Total testcase Written / Total Execution /
# Total Review counts for that block (rows 51-53) now have real values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C51").Value = 1070
$ws.Range("C52").Value = 1243
$ws.Range("C53").Value = 704

# Move the on-screen selection down to the cell that was being worked on
# (reflects the scrolled view/selection captured when the file was saved).
$ws.Range("C53").Select()
